# opiate/data/qa.xlsx -- "update qa and matrix config files"
#
# The ion_ratio_low/ion_ratio_high columns (C/D) are replaced with
# ion_ratio_average/ion_ratio_cv, with new data values throughout.
# K21 (an outlier int_std_peak_area reading) is cleared out entirely.
# Column widths for C/D are adjusted, and D2 becomes the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: rename the ion-ratio columns ---
$ws.Range("C1").Value = "ion_ratio_average"
$ws.Range("D1").Value = "ion_ratio_cv"

# --- new data for columns C (ion_ratio_average) and D (ion_ratio_cv) ---
$ionRatioAverage = @{
    2  = 1.2400000000000002
    3  = 0.40500000000000003
    4  = 1.8800000000000001
    5  = 2.1100000000000003
    6  = 0.28999999999999998
    7  = 3.12
    8  = 1.7200000000000002
    9  = 17.22
    10 = 3.67
    11 = 1.28
    12 = 1.1599999999999999
    13 = 0.97
    14 = 1.5
    15 = 0.2
    16 = 16.329999999999998
    17 = 11.260000000000002
    18 = 3.59
    19 = 3.42
    20 = 11.6
    21 = 10.530000000000001
}

$ionRatioCv = @{
    2  = 0.109
    3  = 0.13900000000000001
    4  = 0.158
    5  = 0.17
    6  = 0.222
    7  = 0.20799999999999999
    8  = 0.072
    9  = 0.16600000000000001
    10 = 0.06
    11 = 0.161
    12 = 0.39300000000000002
    13 = 0.26400000000000001
    14 = 0.14299999999999999
    15 = 0.126
    16 = 0.88
    17 = 0.26
    18 = 0.27300000000000002
    19 = 0.35699999999999998
    20 = 0.41399999999999998
    21 = 0.28199999999999997
}

foreach ($row in 2..21) {
    $ws.Cells.Item($row, 3).Value = $ionRatioAverage[$row]
    $ws.Cells.Item($row, 4).Value = $ionRatioCv[$row]
}

# --- drop the stray int_std_peak_area reading on row 21 ---
$ws.Range("K21").Clear()

# --- column widths: C widens to fit the longer header, D narrows ---
$ws.Columns.Item(3).ColumnWidth = 14.166666666666666
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666

# --- leave the selection on D2, matching the saved view state ---
$ws.Range("D2").Select() | Out-Null
